$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = -1560.156331200843
$ws.Range("C3").Value = 5090.658363836415
$ws.Range("D3").Value = 398.5801982374671
$ws.Range("F3").Value = 3182.312662401685
$ws.Range("G3").Value = 3299.706481224053

$ws.Range("B4").Value = -1509.623188321052
$ws.Range("C4").Value = 1601.071685469473
$ws.Range("D4").Value = 352.0749029812595
$ws.Range("F4").Value = 3113.246376642103
$ws.Range("G4").Value = 3291.230553566339

$ws.Range("B5").Value = -1497.638595834378
$ws.Range("C5").Value = 2269.87698263957
$ws.Range("D5").Value = 349.3206471609823
$ws.Range("F5").Value = 3121.277191668756
$ws.Range("G5").Value = 3359.851726694859

$ws.Range("B6").Value = -1467.353277088631
$ws.Range("C6").Value = 1348.781708161708
$ws.Range("D6").Value = 309.9784430938314
$ws.Range("F6").Value = 3092.706554177263
$ws.Range("G6").Value = 3391.871447305232
